# formatação das anotações com espaçamento 1.5 entre linhas
#
# Apply 1.5 line spacing (wdLineSpace1pt5) to every paragraph in the
# document. In OOXML terms this sets <w:spacing w:line="360"
# w:lineRule="auto"/> on each paragraph's <w:pPr> (merging with any
# existing <w:spacing w:after="0"/> that is already present).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 1   # wdLineSpace1pt5
}
